# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.422.07'
$ws.Range("E2").Value = '  +4.13%  '
$ws.Range("D3").Value = '2.434.93'
$ws.Range("E3").Value = '  +3.20%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.10'
$ws.Range("E5").Value = '  +2.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.64'
$ws.Range("E6").Value = '  +3.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +3.31%  '
$ws.Range("E9").Value = '  +4.85%  '
$ws.Range("E10").Value = '  +3.93%  '
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("E12").Value = '  -2.24%  '
$ws.Range("E13").Value = '  +4.95%  '
$ws.Range("D14").Value = '2.864.89'
$ws.Range("E14").Value = '  +3.12%  '
$ws.Range("D15").Value = '60.323.46'
$ws.Range("E15").Value = '  +4.02%  '
$ws.Range("E16").Value = '  +4.16%  '
$ws.Range("D17").Value = '2.428.49'
$ws.Range("E17").Value = '  +3.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.43'
$ws.Range("E18").Value = '  +5.78%  '
$ws.Range("E19").Value = '  +3.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '335.07'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.83'
$ws.Range("E21").Value = '  +1.74%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.47'
$ws.Range("E23").Value = '  +4.37%  '
$ws.Range("E24").Value = '  +3.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.62'
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").Value = '0.0₃0790'
$ws.Range("E28").Value = '  +6.62%  '
$ws.Range("E29").Value = '  +2.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.34'
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.23'
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.81'
$ws.Range("E32").Value = '  +2.14%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("E35").Value = '  +6.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.24'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '39.85'
$ws.Range("E39").Value = '  +1.00%  '
$ws.Range("E40").Value = '  +10.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '322.13'
$ws.Range("E41").Value = '  +11.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.72'
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.22'
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0962'
$ws.Range("E44").Value = '  +1.73%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0526'
$ws.Range("E45").Value = '  +3.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.64'
$ws.Range("E46").Value = '  +2.40%  '
$ws.Range("E47").Value = '  +8.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.575'
$ws.Range("E48").Value = '  +1.43%  '
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.00'
$ws.Range("E50").Value = '  +2.79%  '
$ws.Range("E51").Value = '  -0.24%  '
